$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 2..18 down to 3..19 (copy including formatting), working
# bottom-up so we don't overwrite rows before they're copied. This avoids the extra
# unused cell style that a native Rows.Insert() would introduce.
for ($r = 18; $r -ge 2; $r--) {
    $src = $ws.Range("A$r" + ":E$r")
    $dst = $ws.Range("A$($r + 1)" + ":E$($r + 1)")
    $src.Copy($dst)
}

# Write the refreshed/recomputed values for every data row (2..19), including the
# newly inserted 2007 forecast row at the top.
$data = @(
    @(2, 39400, 2007, 2.070003986395053, 2008, 0.6967455006573253),
    @(3, 39765, 2008, 0.517569958955022, 2009, -0.6367039903686034),
    @(4, 40130, 2009, -3.956152295564896, 2010, -0.6950853968889392),
    @(5, 40494, 2010, 1.234995474941392, 2011, -0.2098161877568061),
    @(6, 40862, 2011, 0.899360810820804, 2012, 1.205741443109987),
    @(7, 41228, 2012, 0.9010266119894084, 2013, 1.531699207045123),
    @(8, 41592, 2013, 0.02019328874804938, 2014, -0.1259279434590921),
    @(9, 41957, 2014, 0.1729981757035093, 2015, 0.1749537368921361),
    @(10, 42321, 2015, 0.09752710595589686, 2016, -0.001769149545471915),
    @(11, 42689, 2016, -0.5280591151586633, 2017, -0.05116199209030947),
    @(12, 43053, 2017, 0.07201851318385799, 2018, 0.2735900898381383),
    @(13, 43418, 2018, 0.3727661260635617, 2019, -0.9505847809128332),
    @(14, 43783, 2019, -0.801759526476209, 2020, 0.047674034857903),
    @(15, 44159, 2020, -1.103489789942047, 2021, 1.605918384453009),
    @(16, 44525, 2021, 0.9704846793491928, 2022, -0.8255212498362474),
    @(17, 44890, 2022, -0.7009264669202708, 2023, 0.6624163082313173),
    @(18, 45254, 2023, 0.3928252664241905, 2024, 0.302295480375836),
    @(19, 45618, 2024, 0.3224026462283813, 2025, -0.7618983399156787)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
}
